$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "ddd" -> "gg"
$ws.Range("A2").Value = "gg"

# F2: "✗" -> "✔"
$ws.Range("F2").Value = "✔"

# K2: "Приїду на машині" -> "" (cleared)
$ws.Range("K2").Value = ""
